$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 168291
$ws.Range("C4").Value = 159159
$ws.Range("C5").Value = 9132
$ws.Range("C8").Value = 65.5
